# Applies the cryptos.xlsx data refresh described in the commit diff.
# Each (cell, new value) pair below mirrors one changed <c> in the OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text storage for price cells that would
# otherwise be auto-coerced to a number by Excel's type inference
# (e.g. '1.001' -> 1.001 float), matching the source data's string type.
$q = "'"

$ws.Range('D2').Value = '30.789.49'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.938.49'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = $q + '242.12'
$ws.Range('D6').Value = $q + '1.001'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').Value = $q + '0.4882'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').Value = $q + '0.2932'
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('D9').Value = $q + '0.06894'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').Value = $q + '19.56'
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').Value = $q + '105.47'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '1.938.10'
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('D13').Value = $q + '0.07749'
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = $q + '5.316'
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = $q + '0.6973'
$ws.Range('E15').Value = '  -2.82%  '
$ws.Range('D16').Value = $q + '275.56'
$ws.Range('E16').Value = '  -4.43%  '
$ws.Range('D17').Value = '30.792.50'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = $q + '0.000007715'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = $q + '13.13'
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.210.81'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = $q + '1.001'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = $q + '5.436'
$ws.Range('E22').Value = '  -3.80%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = $q + '6.470'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('D25').Value = $q + '9.694'
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('D26').Value = $q + '167.84'
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').Value = $q + '19.58'
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('D29').Value = $q + '0.1041'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('D30').Value = $q + '1.388'
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('E31').Value = '  -3.02%  '
$ws.Range('D32').Value = $q + '4.538'
$ws.Range('D33').Value = $q + '4.355'
$ws.Range('E33').Value = '  -3.82%  '
$ws.Range('D34').Value = $q + '0.04852'
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('D35').Value = $q + '0.7475'
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('D37').Value = $q + '2.726'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('E39').Value = '  -2.02%  '
$ws.Range('D40').Value = $q + '77.32'
$ws.Range('E40').Value = '  +4.49%  '
$ws.Range('D41').Value = $q + '6.419'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = $q + '2.093'
$ws.Range('E42').Value = '  -2.76%  '
$ws.Range('D43').Value = $q + '0.8992'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('D44').Value = $q + '108.12'
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('D45').Value = $q + '0.4404'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').Value = $q + '0.9975'
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').Value = $q + '7.729'
$ws.Range('E47').Value = '  +2.46%  '
$ws.Range('D48').Value = $q + '995.31'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = $q + '9.174'
$ws.Range('E50').Value = '  -2.82%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = $q + '35.67'
$ws.Range('E51').Value = '  -1.37%  '
